$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from serial date 45189 (2023-09-20) to 45190 (2023-09-21)
$ws.Range("C2:C18").Value = 45190
